$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 0.9939999999999999
$ws.Range("E2").Value = 0.006000000000000083
$ws.Range("F2").Value = 0.1099352836608887
$ws.Range("D3").Value = 0.9889999999999998
$ws.Range("E3").Value = 0.01100000000000015
$ws.Range("F3").Value = 0.1381580829620361
$ws.Range("D4").Value = 0.9939999999999999
$ws.Range("E4").Value = 0.006000000000000083
$ws.Range("F4").Value = 0.1155011653900146
$ws.Range("D5").Value = 0.9839999999999999
$ws.Range("E5").Value = 0.01600000000000022
$ws.Range("F5").Value = 0.1315183639526367
$ws.Range("D6").Value = 0.9799999999999998
$ws.Range("E6").Value = 0.02000000000000027
$ws.Range("F6").Value = 0.1159365177154541
$ws.Range("D7").Value = 0.9849999999999999
$ws.Range("E7").Value = 0.0150000000000002
$ws.Range("F7").Value = 0.1224832534790039
$ws.Range("D8").Value = 0.9859999999999999
$ws.Range("E8").Value = 0.01400000000000019
$ws.Range("F8").Value = 0.130143404006958
$ws.Range("D9").Value = 0.9849999999999999
$ws.Range("E9").Value = 0.0150000000000002
$ws.Range("F9").Value = 0.1325314044952393
$ws.Range("D10").Value = 0.9879999999999998
$ws.Range("E10").Value = 0.01200000000000016
$ws.Range("F10").Value = 0.1159682273864746
$ws.Range("D11").Value = 0.9899999999999998
$ws.Range("E11").Value = 0.01000000000000014
$ws.Range("F11").Value = 0.1779632568359375
$ws.Range("D12").Value = 0.9889999999999998
$ws.Range("E12").Value = 0.01100000000000015
$ws.Range("F12").Value = 0.1296448707580566
$ws.Range("D13").Value = 0.9809999999999998
$ws.Range("E13").Value = 0.01900000000000025
$ws.Range("F13").Value = 0.1289284229278564
$ws.Range("D14").Value = 0.9829999999999997
$ws.Range("E14").Value = 0.01700000000000023
$ws.Range("F14").Value = 0.1212136745452881
$ws.Range("D15").Value = 0.9869999999999999
$ws.Range("E15").Value = 0.01300000000000018
$ws.Range("F15").Value = 0.1339349746704102
$ws.Range("D16").Value = 0.9889999999999998
$ws.Range("E16").Value = 0.01100000000000015
$ws.Range("F16").Value = 0.1273660659790039
$ws.Range("D17").Value = 0.9839999999999999
$ws.Range("E17").Value = 0.01600000000000022
$ws.Range("F17").Value = 0.1233963966369629
$ws.Range("D18").Value = 0.9849999999999999
$ws.Range("E18").Value = 0.0150000000000002
$ws.Range("F18").Value = 0.1212573051452637
$ws.Range("D19").Value = 0.992
$ws.Range("E19").Value = 0.008000000000000109
$ws.Range("F19").Value = 0.1214661598205566
$ws.Range("D20").Value = 0.9889999999999998
$ws.Range("E20").Value = 0.01100000000000015
$ws.Range("F20").Value = 0.1505980491638184
$ws.Range("D21").Value = 0.9829999999999997
$ws.Range("E21").Value = 0.01700000000000023
$ws.Range("F21").Value = 0.1359765529632568
$ws.Range("D22").Value = 0.9869999999999999
$ws.Range("E22").Value = 0.01300000000000018
$ws.Range("F22").Value = 0.1159391403198242
$ws.Range("D23").Value = 0.9839999999999999
$ws.Range("E23").Value = 0.01600000000000022
$ws.Range("F23").Value = 0.1315469741821289
$ws.Range("D24").Value = 0.9839999999999999
$ws.Range("E24").Value = 0.01600000000000022
$ws.Range("F24").Value = 0.1205320358276367
$ws.Range("D25").Value = 0.9829999999999997
$ws.Range("E25").Value = 0.01700000000000023
$ws.Range("F25").Value = 0.1311430931091309
$ws.Range("D26").Value = 0.9889999999999998
$ws.Range("E26").Value = 0.01100000000000015
$ws.Range("F26").Value = 0.1213092803955078
$ws.Range("D27").Value = 0.5390832522033016
$ws.Range("E27").Value = 0.4609167477966984
$ws.Range("F27").Value = 0.1293210983276367
$ws.Range("D28").Value = 0.5540769795843318
$ws.Range("E28").Value = 0.4459230204156682
$ws.Range("F28").Value = 0.1215567588806152
$ws.Range("D29").Value = 0.5556899996898556
$ws.Range("E29").Value = 0.4443100003101445
$ws.Range("F29").Value = 0.1245367527008057
$ws.Range("D30").Value = 0.5458365558560589
$ws.Range("E30").Value = 0.454163444143941
$ws.Range("F30").Value = 0.1315212249755859
$ws.Range("D31").Value = 0.5579943471139607
$ws.Range("E31").Value = 0.4420056528860392
$ws.Range("F31").Value = 0.1212725639343262
$ws.Range("D32").Value = 0.5284731420088208
$ws.Range("E32").Value = 0.4715268579911793
$ws.Range("F32").Value = 0.1293597221374512
$ws.Range("D33").Value = 0.5336666111959566
$ws.Range("E33").Value = 0.4663333888040434
$ws.Range("F33").Value = 0.1234517097473145
$ws.Range("D34").Value = 0.5280786779976891
$ws.Range("E34").Value = 0.4719213220023108
$ws.Range("F34").Value = 0.1213588714599609
$ws.Range("D35").Value = 0.5577177341349445
$ws.Range("E35").Value = 0.4422822658650555
$ws.Range("F35").Value = 0.1314828395843506
$ws.Range("D36").Value = 0.5442885979579128
$ws.Range("E36").Value = 0.4557114020420871
$ws.Range("F36").Value = 0.1194169521331787
$ws.Range("D37").Value = 0.5350859687039796
$ws.Range("E37").Value = 0.4649140312960204
$ws.Range("F37").Value = 0.1246562004089355
$ws.Range("D38").Value = 0.5601388396262077
$ws.Range("E38").Value = 0.4398611603737924
$ws.Range("F38").Value = 0.1513404846191406
$ws.Range("D39").Value = 0.5393966168861126
$ws.Range("E39").Value = 0.4606033831138873
$ws.Range("F39").Value = 0.1314418315887451
$ws.Range("D40").Value = 0.5224573373908049
$ws.Range("E40").Value = 0.4775426626091951
$ws.Range("F40").Value = 0.1212208271026611
$ws.Range("D41").Value = 0.5327133433223908
$ws.Range("E41").Value = 0.4672866566776093
$ws.Range("F41").Value = 0.1323482990264893
$ws.Range("D42").Value = 0.5498307296160845
$ws.Range("E42").Value = 0.4501692703839156
$ws.Range("F42").Value = 0.1214070320129395
$ws.Range("D43").Value = 0.5914100124100504
$ws.Range("E43").Value = 0.4085899875899495
$ws.Range("F43").Value = 0.1294624805450439
$ws.Range("D44").Value = 0.5631826025463965
$ws.Range("E44").Value = 0.4368173974536035
$ws.Range("F44").Value = 0.1234138011932373
$ws.Range("D45").Value = 0.5463494380795312
$ws.Range("E45").Value = 0.4536505619204688
$ws.Range("F45").Value = 0.1270678043365479
$ws.Range("D46").Value = 0.5556857233724523
$ws.Range("E46").Value = 0.4443142766275477
$ws.Range("F46").Value = 0.1164133548736572
$ws.Range("D47").Value = 0.5336142616301067
$ws.Range("E47").Value = 0.4663857383698934
$ws.Range("F47").Value = 0.121455192565918
$ws.Range("D48").Value = 0.5680414036803333
$ws.Range("E48").Value = 0.4319585963196667
$ws.Range("F48").Value = 0.1314628124237061
$ws.Range("D49").Value = 0.5585506804142504
$ws.Range("E49").Value = 0.4414493195857497
$ws.Range("F49").Value = 0.1214718818664551
$ws.Range("D50").Value = 0.5617954729119342
$ws.Range("E50").Value = 0.4382045270880658
$ws.Range("F50").Value = 0.132035493850708
$ws.Range("D51").Value = 0.5442440261373676
$ws.Range("E51").Value = 0.4557559738626326
$ws.Range("F51").Value = 0.1211438179016113
$ws.Range("D52").Value = 0.9370000000000001
$ws.Range("E52").Value = 0.063
$ws.Range("F52").Value = 0.2302999496459961
$ws.Range("D53").Value = 0.998
$ws.Range("E53").Value = 0.002
$ws.Range("F53").Value = 0.2148532867431641
$ws.Range("D54").Value = 0.958
$ws.Range("E54").Value = 0.042
$ws.Range("F54").Value = 0.2223718166351318
$ws.Range("D55").Value = 1
$ws.Range("E55").Value = 0
$ws.Range("F55").Value = 0.2224123477935791
$ws.Range("D56").Value = 0.979
$ws.Range("E56").Value = 0.021
$ws.Range("F56").Value = 0.2206299304962158
$ws.Range("D57").Value = 0.978
$ws.Range("E57").Value = 0.022
$ws.Range("F57").Value = 0.2143073081970215
$ws.Range("D58").Value = 0.998
$ws.Range("E58").Value = 0.002
$ws.Range("F58").Value = 0.2131123542785645
$ws.Range("D59").Value = 0.997
$ws.Range("E59").Value = 0.003
$ws.Range("F59").Value = 0.2225186824798584
$ws.Range("D60").Value = 1
$ws.Range("E60").Value = 0
$ws.Range("F60").Value = 0.2127552032470703
$ws.Range("D61").Value = 0.997
$ws.Range("E61").Value = 0.003
$ws.Range("F61").Value = 0.2228224277496338
$ws.Range("D62").Value = 1
$ws.Range("E62").Value = 0
$ws.Range("F62").Value = 0.2202463150024414
$ws.Range("D63").Value = 0.977
$ws.Range("E63").Value = 0.023
$ws.Range("F63").Value = 0.2149455547332764
$ws.Range("D64").Value = 0.981
$ws.Range("E64").Value = 0.019
$ws.Range("F64").Value = 0.2123098373413086
$ws.Range("D65").Value = 0.99
$ws.Range("E65").Value = 0.01
$ws.Range("F65").Value = 0.220484733581543
$ws.Range("D66").Value = 0.99
$ws.Range("E66").Value = 0.01
$ws.Range("F66").Value = 0.2144982814788818
$ws.Range("D67").Value = 0.981
$ws.Range("E67").Value = 0.019
$ws.Range("F67").Value = 0.2223618030548096
$ws.Range("D68").Value = 0.982
$ws.Range("E68").Value = 0.018
$ws.Range("F68").Value = 0.2224512100219727
$ws.Range("D69").Value = 0.98
$ws.Range("E69").Value = 0.02
$ws.Range("F69").Value = 0.2123129367828369
$ws.Range("D70").Value = 0.981
$ws.Range("E70").Value = 0.019
$ws.Range("F70").Value = 0.2154088020324707
$ws.Range("D71").Value = 1
$ws.Range("E71").Value = 0
$ws.Range("F71").Value = 0.2160854339599609
$ws.Range("D72").Value = 0.974
$ws.Range("E72").Value = 0.026
$ws.Range("F72").Value = 0.2142362594604492
$ws.Range("D73").Value = 0.986
$ws.Range("E73").Value = 0.014
$ws.Range("F73").Value = 0.2181284427642822
$ws.Range("D74").Value = 0.948
$ws.Range("E74").Value = 0.052
$ws.Range("F74").Value = 0.2160434722900391
$ws.Range("D75").Value = 0.983
$ws.Range("E75").Value = 0.017
$ws.Range("F75").Value = 0.2226972579956055
$ws.Range("D76").Value = 1
$ws.Range("E76").Value = 0
$ws.Range("F76").Value = 0.2216298580169678
$ws.Range("D77").Value = 0.268
$ws.Range("E77").Value = 0.732
$ws.Range("F77").Value = 0.225304126739502
$ws.Range("D78").Value = 0.535
$ws.Range("E78").Value = 0.465
$ws.Range("F78").Value = 0.2325074672698975
$ws.Range("D79").Value = 0.622
$ws.Range("E79").Value = 0.378
$ws.Range("F79").Value = 0.2236850261688232
$ws.Range("D80").Value = 0.705
$ws.Range("E80").Value = 0.295
$ws.Range("F80").Value = 0.222783088684082
$ws.Range("D81").Value = 0.654
$ws.Range("E81").Value = 0.346
$ws.Range("F81").Value = 0.2204570770263672
$ws.Range("D82").Value = 0.628
$ws.Range("E82").Value = 0.372
$ws.Range("F82").Value = 0.2157690525054932
$ws.Range("D83").Value = 0.433
$ws.Range("E83").Value = 0.5669999999999999
$ws.Range("F83").Value = 0.2224600315093994
$ws.Range("D84").Value = 0.313
$ws.Range("E84").Value = 0.6870000000000001
$ws.Range("F84").Value = 0.2233648300170898
$ws.Range("D85").Value = 0.852
$ws.Range("E85").Value = 0.148
$ws.Range("F85").Value = 0.2224650382995605
$ws.Range("D86").Value = 0.595
$ws.Range("E86").Value = 0.405
$ws.Range("F86").Value = 0.2323131561279297
$ws.Range("D87").Value = 0.425
$ws.Range("E87").Value = 0.575
$ws.Range("F87").Value = 0.222618579864502
$ws.Range("D88").Value = 0.503
$ws.Range("E88").Value = 0.497
$ws.Range("F88").Value = 0.2325046062469482
$ws.Range("D89").Value = 0.879
$ws.Range("E89").Value = 0.121
$ws.Range("F89").Value = 0.2237091064453125
$ws.Range("D90").Value = 0.5639999999999999
$ws.Range("E90").Value = 0.436
$ws.Range("F90").Value = 0.2222957611083984
$ws.Range("D91").Value = 0.573
$ws.Range("E91").Value = 0.427
$ws.Range("F91").Value = 0.2292194366455078
$ws.Range("D92").Value = 0.834
$ws.Range("E92").Value = 0.166
$ws.Range("F92").Value = 0.2227916717529297
$ws.Range("D93").Value = 0.8120000000000001
$ws.Range("E93").Value = 0.188
$ws.Range("F93").Value = 0.2272853851318359
$ws.Range("D94").Value = 0.252
$ws.Range("E94").Value = 0.748
$ws.Range("F94").Value = 0.2226216793060303
$ws.Range("D95").Value = 0.532
$ws.Range("E95").Value = 0.468
$ws.Range("F95").Value = 0.2161240577697754
$ws.Range("D96").Value = 0.362
$ws.Range("E96").Value = 0.638
$ws.Range("F96").Value = 0.2159979343414307
$ws.Range("D97").Value = 0.647
$ws.Range("E97").Value = 0.353
$ws.Range("F97").Value = 0.2264444828033447
$ws.Range("D98").Value = 0.433
$ws.Range("E98").Value = 0.5669999999999999
$ws.Range("F98").Value = 0.2212343215942383
$ws.Range("D99").Value = 0.786
$ws.Range("E99").Value = 0.214
$ws.Range("F99").Value = 0.2161719799041748
$ws.Range("D100").Value = 0.401
$ws.Range("E100").Value = 0.599
$ws.Range("F100").Value = 0.2161731719970703
$ws.Range("D101").Value = 0.721
$ws.Range("E101").Value = 0.279
$ws.Range("F101").Value = 0.2338113784790039
$ws.Range("D102").Value = 0.986013986013986
$ws.Range("E102").Value = 0.01398601398601399
$ws.Range("F102").Value = 42.34375
$ws.Range("D103").Value = 0.983016983016983
$ws.Range("E103").Value = 0.01698301698301698
$ws.Range("F103").Value = 42.8125
$ws.Range("D104").Value = 0.986013986013986
$ws.Range("E104").Value = 0.01398601398601399
$ws.Range("F104").Value = 43.28125
$ws.Range("D105").Value = 0.9890109890109891
$ws.Range("E105").Value = 0.01098901098901099
$ws.Range("F105").Value = 41.71875
$ws.Range("D106").Value = 0.987012987012987
$ws.Range("E106").Value = 0.01298701298701299
$ws.Range("F106").Value = 43.8125
$ws.Range("D107").Value = 0.99000999000999
$ws.Range("E107").Value = 0.00999000999000999
$ws.Range("F107").Value = 42.65625
$ws.Range("D108").Value = 0.988011988011988
$ws.Range("E108").Value = 0.01198801198801199
$ws.Range("F108").Value = 44.171875
$ws.Range("D109").Value = 0.967032967032967
$ws.Range("E109").Value = 0.03296703296703297
$ws.Range("F109").Value = 42.375
$ws.Range("D110").Value = 0.9820179820179821
$ws.Range("E110").Value = 0.01798201798201798
$ws.Range("F110").Value = 42.96875
$ws.Range("D111").Value = 0.983016983016983
$ws.Range("E111").Value = 0.01698301698301698
$ws.Range("F111").Value = 42.75
$ws.Range("D112").Value = 0.991008991008991
$ws.Range("E112").Value = 0.008991008991008992
$ws.Range("F112").Value = 44.109375
$ws.Range("D113").Value = 0.981018981018981
$ws.Range("E113").Value = 0.01898101898101898
$ws.Range("F113").Value = 44.859375
$ws.Range("D114").Value = 0.9920079920079921
$ws.Range("E114").Value = 0.007992007992007992
$ws.Range("F114").Value = 42.40625
$ws.Range("D115").Value = 0.981018981018981
$ws.Range("E115").Value = 0.01898101898101898
$ws.Range("F115").Value = 44.359375
$ws.Range("D116").Value = 0.991008991008991
$ws.Range("E116").Value = 0.008991008991008992
$ws.Range("F116").Value = 42.0625
$ws.Range("D117").Value = 0.991008991008991
$ws.Range("E117").Value = 0.008991008991008992
$ws.Range("F117").Value = 43.625
$ws.Range("D118").Value = 0.994005994005994
$ws.Range("E118").Value = 0.005994005994005994
$ws.Range("F118").Value = 41.09375
$ws.Range("D119").Value = 0.977022977022977
$ws.Range("E119").Value = 0.02297702297702298
$ws.Range("F119").Value = 42.5625
$ws.Range("D120").Value = 0.985014985014985
$ws.Range("E120").Value = 0.01498501498501499
$ws.Range("F120").Value = 43.109375
$ws.Range("D121").Value = 0.988011988011988
$ws.Range("E121").Value = 0.01198801198801199
$ws.Range("F121").Value = 41.765625
$ws.Range("D122").Value = 0.993006993006993
$ws.Range("E122").Value = 0.006993006993006993
$ws.Range("F122").Value = 41.390625
$ws.Range("D123").Value = 0.985014985014985
$ws.Range("E123").Value = 0.01498501498501499
$ws.Range("F123").Value = 41.953125
$ws.Range("D124").Value = 0.987012987012987
$ws.Range("E124").Value = 0.01298701298701299
$ws.Range("F124").Value = 42.4375
$ws.Range("D125").Value = 0.984015984015984
$ws.Range("E125").Value = 0.01598401598401598
$ws.Range("F125").Value = 41.796875
$ws.Range("D126").Value = 0.9890109890109891
$ws.Range("E126").Value = 0.01098901098901099
$ws.Range("F126").Value = 40.328125
$ws.Range("D127").Value = 0.6123876123876124
$ws.Range("E127").Value = 0.3876123876123876
$ws.Range("F127").Value = 43.453125
$ws.Range("D128").Value = 0.6243756243756243
$ws.Range("E128").Value = 0.3756243756243756
$ws.Range("F128").Value = 44.890625
$ws.Range("D129").Value = 0.6543456543456544
$ws.Range("E129").Value = 0.3456543456543457
$ws.Range("F129").Value = 43.96875
$ws.Range("D130").Value = 0.6673326673326674
$ws.Range("E130").Value = 0.3326673326673327
$ws.Range("F130").Value = 44.78125
$ws.Range("D131").Value = 0.6253746253746254
$ws.Range("E131").Value = 0.3746253746253747
$ws.Range("F131").Value = 42.625
$ws.Range("D132").Value = 0.6103896103896104
$ws.Range("E132").Value = 0.3896103896103896
$ws.Range("F132").Value = 43.59375
$ws.Range("D133").Value = 0.6483516483516484
$ws.Range("E133").Value = 0.3516483516483517
$ws.Range("F133").Value = 44.84375
$ws.Range("D134").Value = 0.6213786213786214
$ws.Range("E134").Value = 0.3786213786213786
$ws.Range("F134").Value = 43.59375
$ws.Range("D135").Value = 0.6263736263736264
$ws.Range("E135").Value = 0.3736263736263736
$ws.Range("F135").Value = 45.28125
$ws.Range("D136").Value = 0.6223776223776224
$ws.Range("E136").Value = 0.3776223776223776
$ws.Range("F136").Value = 44.046875
$ws.Range("D137").Value = 0.6243756243756243
$ws.Range("E137").Value = 0.3756243756243756
$ws.Range("F137").Value = 44.515625
$ws.Range("D138").Value = 0.6333666333666333
$ws.Range("E138").Value = 0.3666333666333667
$ws.Range("F138").Value = 44.65625
$ws.Range("D139").Value = 0.6143856143856143
$ws.Range("E139").Value = 0.3856143856143856
$ws.Range("F139").Value = 43.1875
$ws.Range("D140").Value = 0.6533466533466533
$ws.Range("E140").Value = 0.3466533466533466
$ws.Range("F140").Value = 45.4375
$ws.Range("D141").Value = 0.6213786213786214
$ws.Range("E141").Value = 0.3786213786213786
$ws.Range("F141").Value = 45.765625
$ws.Range("D142").Value = 0.6313686313686314
$ws.Range("E142").Value = 0.3686313686313686
$ws.Range("F142").Value = 43.265625
$ws.Range("D143").Value = 0.6593406593406593
$ws.Range("E143").Value = 0.3406593406593407
$ws.Range("F143").Value = 43.890625
$ws.Range("D144").Value = 0.6263736263736264
$ws.Range("E144").Value = 0.3736263736263736
$ws.Range("F144").Value = 44.671875
$ws.Range("D145").Value = 0.6393606393606394
$ws.Range("E145").Value = 0.3606393606393606
$ws.Range("F145").Value = 42.78125
$ws.Range("D146").Value = 0.6063936063936064
$ws.Range("E146").Value = 0.3936063936063936
$ws.Range("F146").Value = 43.6875
$ws.Range("D147").Value = 0.6373626373626373
$ws.Range("E147").Value = 0.3626373626373626
$ws.Range("F147").Value = 44.65625
$ws.Range("D148").Value = 0.6463536463536463
$ws.Range("E148").Value = 0.3536463536463537
$ws.Range("F148").Value = 44.71875
$ws.Range("D149").Value = 0.6373626373626373
$ws.Range("E149").Value = 0.3626373626373626
$ws.Range("F149").Value = 44.578125
$ws.Range("D150").Value = 0.6393606393606394
$ws.Range("E150").Value = 0.3606393606393606
$ws.Range("F150").Value = 43.109375
$ws.Range("D151").Value = 0.6433566433566433
$ws.Range("E151").Value = 0.3566433566433567
$ws.Range("F151").Value = 42.375
$ws.Range("D152").Value = 0.986013986013986
$ws.Range("E152").Value = 0.01398601398601399
$ws.Range("F152").Value = 49.28125
$ws.Range("D153").Value = 0.984015984015984
$ws.Range("E153").Value = 0.01598401598401598
$ws.Range("F153").Value = 48.234375
$ws.Range("D154").Value = 0.991008991008991
$ws.Range("E154").Value = 0.008991008991008992
$ws.Range("F154").Value = 48.53125
$ws.Range("D155").Value = 0.986013986013986
$ws.Range("E155").Value = 0.01398601398601399
$ws.Range("F155").Value = 48.125
$ws.Range("D156").Value = 0.983016983016983
$ws.Range("E156").Value = 0.01698301698301698
$ws.Range("F156").Value = 47.546875
$ws.Range("D157").Value = 0.987012987012987
$ws.Range("E157").Value = 0.01298701298701299
$ws.Range("F157").Value = 48.21875
$ws.Range("D158").Value = 0.987012987012987
$ws.Range("E158").Value = 0.01298701298701299
$ws.Range("F158").Value = 48.5625
$ws.Range("D159").Value = 0.995004995004995
$ws.Range("E159").Value = 0.004995004995004995
$ws.Range("F159").Value = 47.609375
$ws.Range("D160").Value = 0.991008991008991
$ws.Range("E160").Value = 0.008991008991008992
$ws.Range("F160").Value = 48.5625
$ws.Range("D161").Value = 0.987012987012987
$ws.Range("E161").Value = 0.01298701298701299
$ws.Range("F161").Value = 48.21875
$ws.Range("D162").Value = 0.987012987012987
$ws.Range("E162").Value = 0.01298701298701299
$ws.Range("F162").Value = 49.96875
$ws.Range("D163").Value = 0.987012987012987
$ws.Range("E163").Value = 0.01298701298701299
$ws.Range("F163").Value = 47.921875
$ws.Range("D164").Value = 0.986013986013986
$ws.Range("E164").Value = 0.01398601398601399
$ws.Range("F164").Value = 50.328125
$ws.Range("D165").Value = 0.9890109890109891
$ws.Range("E165").Value = 0.01098901098901099
$ws.Range("F165").Value = 48.796875
$ws.Range("D166").Value = 0.983016983016983
$ws.Range("E166").Value = 0.01698301698301698
$ws.Range("F166").Value = 47.859375
$ws.Range("D167").Value = 0.986013986013986
$ws.Range("E167").Value = 0.01398601398601399
$ws.Range("F167").Value = 48.84375
$ws.Range("D168").Value = 0.986013986013986
$ws.Range("E168").Value = 0.01398601398601399
$ws.Range("F168").Value = 49.046875
$ws.Range("D169").Value = 0.9920079920079921
$ws.Range("E169").Value = 0.007992007992007992
$ws.Range("F169").Value = 49.5625
$ws.Range("D170").Value = 0.986013986013986
$ws.Range("E170").Value = 0.01398601398601399
$ws.Range("F170").Value = 48.078125
$ws.Range("D171").Value = 0.9790209790209791
$ws.Range("E171").Value = 0.02097902097902098
$ws.Range("F171").Value = 48.96875
$ws.Range("D172").Value = 0.987012987012987
$ws.Range("E172").Value = 0.01298701298701299
$ws.Range("F172").Value = 48.96875
$ws.Range("D173").Value = 0.995004995004995
$ws.Range("E173").Value = 0.004995004995004995
$ws.Range("F173").Value = 48.40625
$ws.Range("D174").Value = 0.983016983016983
$ws.Range("E174").Value = 0.01698301698301698
$ws.Range("F174").Value = 49
$ws.Range("D175").Value = 0.983016983016983
$ws.Range("E175").Value = 0.01698301698301698
$ws.Range("F175").Value = 45.53125
$ws.Range("D176").Value = 0.98001998001998
$ws.Range("E176").Value = 0.01998001998001998
$ws.Range("F176").Value = 48.34375
$ws.Range("D177").Value = 0.6273726273726273
$ws.Range("E177").Value = 0.3726273726273726
$ws.Range("F177").Value = 48.59375
$ws.Range("D178").Value = 0.6183816183816184
$ws.Range("E178").Value = 0.3816183816183816
$ws.Range("F178").Value = 48.484375
$ws.Range("D179").Value = 0.6193806193806194
$ws.Range("E179").Value = 0.3806193806193806
$ws.Range("F179").Value = 49.609375
$ws.Range("D180").Value = 0.6483516483516484
$ws.Range("E180").Value = 0.3516483516483517
$ws.Range("F180").Value = 48.984375
$ws.Range("D181").Value = 0.6343656343656343
$ws.Range("E181").Value = 0.3656343656343656
$ws.Range("F181").Value = 49.28125
$ws.Range("D182").Value = 0.6303696303696303
$ws.Range("E182").Value = 0.3696303696303696
$ws.Range("F182").Value = 50.46875
$ws.Range("D183").Value = 0.6483516483516484
$ws.Range("E183").Value = 0.3516483516483517
$ws.Range("F183").Value = 49.125
$ws.Range("D184").Value = 0.6383616383616384
$ws.Range("E184").Value = 0.3616383616383617
$ws.Range("F184").Value = 49.609375
$ws.Range("D185").Value = 0.6303696303696303
$ws.Range("E185").Value = 0.3696303696303696
$ws.Range("F185").Value = 49.390625
$ws.Range("D186").Value = 0.6643356643356644
$ws.Range("E186").Value = 0.3356643356643357
$ws.Range("F186").Value = 48.84375
$ws.Range("D187").Value = 0.6253746253746254
$ws.Range("E187").Value = 0.3746253746253747
$ws.Range("F187").Value = 49.1875
$ws.Range("D188").Value = 0.6423576423576424
$ws.Range("E188").Value = 0.3576423576423576
$ws.Range("F188").Value = 49.25
$ws.Range("D189").Value = 0.6063936063936064
$ws.Range("E189").Value = 0.3936063936063936
$ws.Range("F189").Value = 50.203125
$ws.Range("D190").Value = 0.6423576423576424
$ws.Range("E190").Value = 0.3576423576423576
$ws.Range("F190").Value = 49.890625
$ws.Range("D191").Value = 0.6413586413586414
$ws.Range("E191").Value = 0.3586413586413587
$ws.Range("F191").Value = 50.5
$ws.Range("D192").Value = 0.6593406593406593
$ws.Range("E192").Value = 0.3406593406593407
$ws.Range("F192").Value = 49.9375
$ws.Range("D193").Value = 0.6273726273726273
$ws.Range("E193").Value = 0.3726273726273726
$ws.Range("F193").Value = 49.6875
$ws.Range("D194").Value = 0.5854145854145855
$ws.Range("E194").Value = 0.4145854145854146
$ws.Range("F194").Value = 49.875
$ws.Range("D195").Value = 0.6333666333666333
$ws.Range("E195").Value = 0.3666333666333667
$ws.Range("F195").Value = 49.703125
$ws.Range("D196").Value = 0.6333666333666333
$ws.Range("E196").Value = 0.3666333666333667
$ws.Range("F196").Value = 49.5
$ws.Range("D197").Value = 0.6153846153846154
$ws.Range("E197").Value = 0.3846153846153846
$ws.Range("F197").Value = 49.15625
$ws.Range("D198").Value = 0.6283716283716284
$ws.Range("E198").Value = 0.3716283716283716
$ws.Range("F198").Value = 47.890625
$ws.Range("D199").Value = 0.6013986013986014
$ws.Range("E199").Value = 0.3986013986013986
$ws.Range("F199").Value = 48.703125
$ws.Range("D200").Value = 0.6523476523476524
$ws.Range("E200").Value = 0.3476523476523476
$ws.Range("F200").Value = 49.375
$ws.Range("D201").Value = 0.6293706293706294
$ws.Range("E201").Value = 0.3706293706293706
$ws.Range("F201").Value = 46.84375
$ws.Range("D202").Value = 0.988011988011988
$ws.Range("E202").Value = 0.01198801198801199
$ws.Range("F202").Value = 54.6875
$ws.Range("D203").Value = 0.99000999000999
$ws.Range("E203").Value = 0.00999000999000999
$ws.Range("F203").Value = 54.203125
$ws.Range("D204").Value = 0.988011988011988
$ws.Range("E204").Value = 0.01198801198801199
$ws.Range("F204").Value = 53.875
$ws.Range("D205").Value = 0.991008991008991
$ws.Range("E205").Value = 0.008991008991008992
$ws.Range("F205").Value = 54.15625
$ws.Range("D206").Value = 0.986013986013986
$ws.Range("E206").Value = 0.01398601398601399
$ws.Range("F206").Value = 53.625
$ws.Range("D207").Value = 0.994005994005994
$ws.Range("E207").Value = 0.005994005994005994
$ws.Range("F207").Value = 55.390625
$ws.Range("D208").Value = 0.99000999000999
$ws.Range("E208").Value = 0.00999000999000999
$ws.Range("F208").Value = 54.796875
$ws.Range("D209").Value = 0.9820179820179821
$ws.Range("E209").Value = 0.01798201798201798
$ws.Range("F209").Value = 54.3125
$ws.Range("D210").Value = 0.993006993006993
$ws.Range("E210").Value = 0.006993006993006993
$ws.Range("F210").Value = 54.921875
$ws.Range("D211").Value = 0.993006993006993
$ws.Range("E211").Value = 0.006993006993006993
$ws.Range("F211").Value = 54.65625
$ws.Range("D212").Value = 0.9920079920079921
$ws.Range("E212").Value = 0.007992007992007992
$ws.Range("F212").Value = 54.34375
$ws.Range("D213").Value = 0.988011988011988
$ws.Range("E213").Value = 0.01198801198801199
$ws.Range("F213").Value = 54.28125
$ws.Range("D214").Value = 0.9920079920079921
$ws.Range("E214").Value = 0.007992007992007992
$ws.Range("F214").Value = 53.515625
$ws.Range("D215").Value = 0.9920079920079921
$ws.Range("E215").Value = 0.007992007992007992
$ws.Range("F215").Value = 54.84375
$ws.Range("D216").Value = 0.981018981018981
$ws.Range("E216").Value = 0.01898101898101898
$ws.Range("F216").Value = 53.3125
$ws.Range("D217").Value = 0.9890109890109891
$ws.Range("E217").Value = 0.01098901098901099
$ws.Range("F217").Value = 54.421875
$ws.Range("D218").Value = 0.985014985014985
$ws.Range("E218").Value = 0.01498501498501499
$ws.Range("F218").Value = 54.8125
$ws.Range("D219").Value = 0.987012987012987
$ws.Range("E219").Value = 0.01298701298701299
$ws.Range("F219").Value = 54.765625
$ws.Range("D220").Value = 0.9820179820179821
$ws.Range("E220").Value = 0.01798201798201798
$ws.Range("F220").Value = 55.421875
$ws.Range("D221").Value = 0.984015984015984
$ws.Range("E221").Value = 0.01598401598401598
$ws.Range("F221").Value = 53.765625
$ws.Range("D222").Value = 0.981018981018981
$ws.Range("E222").Value = 0.01898101898101898
$ws.Range("F222").Value = 54.890625
$ws.Range("D223").Value = 0.985014985014985
$ws.Range("E223").Value = 0.01498501498501499
$ws.Range("F223").Value = 54.546875
$ws.Range("D224").Value = 0.99000999000999
$ws.Range("E224").Value = 0.00999000999000999
$ws.Range("F224").Value = 54.96875
$ws.Range("D225").Value = 0.9820179820179821
$ws.Range("E225").Value = 0.01798201798201798
$ws.Range("F225").Value = 52.84375
$ws.Range("D226").Value = 0.98001998001998
$ws.Range("E226").Value = 0.01998001998001998
$ws.Range("F226").Value = 53.59375
$ws.Range("D227").Value = 0.6353646353646354
$ws.Range("E227").Value = 0.3646353646353646
$ws.Range("F227").Value = 56.015625
$ws.Range("D228").Value = 0.6513486513486514
$ws.Range("E228").Value = 0.3486513486513487
$ws.Range("F228").Value = 55.828125
$ws.Range("D229").Value = 0.6463536463536463
$ws.Range("E229").Value = 0.3536463536463537
$ws.Range("F229").Value = 55.875
$ws.Range("D230").Value = 0.6473526473526473
$ws.Range("E230").Value = 0.3526473526473526
$ws.Range("F230").Value = 55.546875
$ws.Range("D231").Value = 0.6363636363636364
$ws.Range("E231").Value = 0.3636363636363636
$ws.Range("F231").Value = 55.234375
$ws.Range("D232").Value = 0.6413586413586414
$ws.Range("E232").Value = 0.3586413586413587
$ws.Range("F232").Value = 56.65625
$ws.Range("D233").Value = 0.6793206793206793
$ws.Range("E233").Value = 0.3206793206793207
$ws.Range("F233").Value = 55.71875
$ws.Range("D234").Value = 0.6313686313686314
$ws.Range("E234").Value = 0.3686313686313686
$ws.Range("F234").Value = 55.3125
$ws.Range("D235").Value = 0.6813186813186813
$ws.Range("E235").Value = 0.3186813186813187
$ws.Range("F235").Value = 55.28125
$ws.Range("D236").Value = 0.6353646353646354
$ws.Range("E236").Value = 0.3646353646353646
$ws.Range("F236").Value = 55.5
$ws.Range("D237").Value = 0.6503496503496503
$ws.Range("E237").Value = 0.3496503496503496
$ws.Range("F237").Value = 55.875
$ws.Range("D238").Value = 0.6673326673326674
$ws.Range("E238").Value = 0.3326673326673327
$ws.Range("F238").Value = 56.15625
$ws.Range("D239").Value = 0.6813186813186813
$ws.Range("E239").Value = 0.3186813186813187
$ws.Range("F239").Value = 55.53125
$ws.Range("D240").Value = 0.6553446553446554
$ws.Range("E240").Value = 0.3446553446553446
$ws.Range("F240").Value = 54.921875
$ws.Range("D241").Value = 0.6433566433566433
$ws.Range("E241").Value = 0.3566433566433567
$ws.Range("F241").Value = 55.328125
$ws.Range("D242").Value = 0.6073926073926074
$ws.Range("E242").Value = 0.3926073926073926
$ws.Range("F242").Value = 56.9375
$ws.Range("D243").Value = 0.6403596403596403
$ws.Range("E243").Value = 0.3596403596403596
$ws.Range("F243").Value = 54.390625
$ws.Range("D244").Value = 0.6293706293706294
$ws.Range("E244").Value = 0.3706293706293706
$ws.Range("F244").Value = 54.359375
$ws.Range("D245").Value = 0.6583416583416584
$ws.Range("E245").Value = 0.3416583416583416
$ws.Range("F245").Value = 55.875
$ws.Range("D246").Value = 0.6413586413586414
$ws.Range("E246").Value = 0.3586413586413587
$ws.Range("F246").Value = 54.53125
$ws.Range("D247").Value = 0.6673326673326674
$ws.Range("E247").Value = 0.3326673326673327
$ws.Range("F247").Value = 54.609375
$ws.Range("D248").Value = 0.6383616383616384
$ws.Range("E248").Value = 0.3616383616383617
$ws.Range("F248").Value = 55.5
$ws.Range("D249").Value = 0.6413586413586414
$ws.Range("E249").Value = 0.3586413586413587
$ws.Range("F249").Value = 54.71875
$ws.Range("D250").Value = 0.6423576423576424
$ws.Range("E250").Value = 0.3576423576423576
$ws.Range("F250").Value = 54.59375
$ws.Range("D251").Value = 0.6553446553446554
$ws.Range("E251").Value = 0.3446553446553446
$ws.Range("F251").Value = 53.6875
